$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the price/volume cells remain stored as plain text (matching source format)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "299.86"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-1.37%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "31.44"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-1.76%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.146"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-2.83%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07315"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-1.82%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.806"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "22.30%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.774"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.82%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.737"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.96%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9255"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.10%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1679"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-1.13%"

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-6.96%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08108"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.61%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03002"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-1.03%"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09908"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.63%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001487"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.74%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006170"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-2.37%"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.458"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.54%"

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.24%"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-2.29%"

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-2.36%"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.563"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.51%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04644"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "1.58%"

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-3.71%"

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.38%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004737"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "7.42%"

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-7.42%"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001873"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "5.34%"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01718"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-1.12%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04497"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-0.42%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007088"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.37%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1335"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-0.21%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002226"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-0.74%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01047"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-22.39%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006221"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-0.06%"

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-21.53%"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.7394"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "4.26%"

